# Slide 3, Shape 1 ("Shape 27" / id=27) - the "What the Algorithm Found" header box.
# The author resized/repositioned the box slightly and replaced the headline text
# (PowerPoint's normAutofit shrank the font afterwards to keep the new, longer
# text on one visual block - fontScale/lnSpcReduction are display-only values
# PowerPoint computes at layout time and are not writable through the object
# model, so they are left for the host's own relayout).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(1)

# Reposition / resize the placeholder (values chosen so the Single-precision
# point round-trip lands exactly on the target EMUs: 2479488,656904 / 9553486,946610).
$shp.Left   = 195.23532104492188
$shp.Top    = 51.72476577758789
$shp.Width  = 752.2430419921875
$shp.Height = 74.53626251220703

# Replace only the headline run's text (keep the trailing empty paragraph intact).
$tr = $shp.TextFrame.TextRange
$old = "What the Algorithm Found"
$headline = $tr.Characters(1, $old.Length)
$headline.Text = "What was discovered in the Data"
